# Apply the "aanpassingen gemaakt aan de logincheck en login class" edit
# to the week49 sheet of the logboek workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week49")

# Row 7: update eindtijd (D7), add activiteiten text (F7) and let the
# duration formula in G7 recalc automatically.
$ws.Cells.Item(7, 4).Value = 0.635416666666667
$ws.Cells.Item(7, 6).Value = "loginclasstest maken met button."

# Row 8: add a new day entry (vrijdag 6 dec 2013)
$ws.Cells.Item(8, 1).Value = "vrijdag"
$ws.Cells.Item(8, 2).Value = 41614
$ws.Cells.Item(8, 2).NumberFormat = "DD/MM/YY"

# The filled-in row now needs the taller "content" row height used
# elsewhere in the workbook for fully populated rows.
$ws.Rows.Item(7).RowHeight = 13.4

# Move the active selection to J3 as in the saved workbook
$ws.Range("J3").Select()

$wb.Save()
